$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Comment text used for the "Comments" column (AC) explaining the scale book corrections.
$msg17377 = "Changed scale book number from 17377 to 1036314 based on results from Age Batch File. Not sure why the difference. Nick Brown Nov 2023."
$msg17376 = "Changed scale book number from 17376 to 1036313 based on results from Age Batch File. Not sure why the difference. Nick Brown Nov 2023."

# Rows 104-113: Scale Book No changed from 17376 to 1036313.
# Write the 17377 message first (on row 114, below) so the shared-string table
# ends up with the 17377 message at index 161 and the 17376 message at index 162,
# matching the order the strings were appended to the workbook.

# Rows 114-123: Scale Book No changed from 17377 to 1036314 (handled first, see above note).
for ($row = 114; $row -le 123; $row++) {
    $ws.Range("E$row").Value = 1036314
    $ws.Range("AC$row").Value = $msg17377
}

for ($row = 104; $row -le 113; $row++) {
    $ws.Range("E$row").Value = 1036313
    $ws.Range("AC$row").Value = $msg17376
}

# Row 110 lost its special (highlighted) row formatting on column E and had its
# column AC formatting reset to the plain bordered style used elsewhere in the
# table when it was edited in Excel. Reproduce that here.
$ws.Range("E110").Style = "Normal"
$ws.Range("AC104").Copy()
$ws.Range("AC110").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update the active selection to reflect where the editor ended up after making
# these edits.
$ws.Range("H124").Select()
